$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-11-27"

# Update the row label for November.
$ws.Range("A13").Value = "November (through 11-27)"

# Row 13 - November (through 11-27)
$ws.Range("C13").Value = 28
$ws.Range("D13").Value = 0.0345
$ws.Range("F13").Value = 63
$ws.Range("G13").Value = 0.1
$ws.Range("I13").Value = 98
$ws.Range("J13").Value = 0.02
$ws.Range("L13").Value = 52
$ws.Range("M13").Value = 0.1333
$ws.Range("O13").Value = 41
$ws.Range("P13").Value = 0.1277
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = 182
$ws.Range("S13").Value = 0.0521
$ws.Range("U13").Value = 180
$ws.Range("V13").Value = 0.0217

# Row 14 - Total
$ws.Range("C14").Value = 254
$ws.Range("D14").Value = 0.115
$ws.Range("F14").Value = 497
$ws.Range("G14").Value = 0.1061
$ws.Range("I14").Value = 747
$ws.Range("J14").Value = 0.0778
$ws.Range("L14").Value = 601
$ws.Range("M14").Value = 0.1096
$ws.Range("O14").Value = 475
$ws.Range("P14").Value = 0.1021
$ws.Range("Q14").Value = 63
$ws.Range("R14").Value = 1186
$ws.Range("S14").Value = 0.0504
$ws.Range("U14").Value = 1531
$ws.Range("V14").Value = 0.059
